$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z1").Borders.Item(7).LineStyle = 1
$ws.Range("Z1").Borders.Item(7).Weight = -4138
$ws.Range("Z1").Borders.Item(9).LineStyle = 1
$ws.Range("Z1").Borders.Item(9).Weight = -4138
$ws.Range("Z1").Value = "test"
